$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (table is sorted ascending by column O, and this
# new CPU's rating puts it right after row 2), pushing the rest of the data
# down by one row.
$ws.Rows(3).Insert()

# Fill in the new benchmark entry: Ryzen 5 3600X
$ws.Range("A3").Value = "AMD"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Ryzen 5 3600X"
$ws.Range("D3").Value = 95
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 3.8
$ws.Range("H3").Value = 4.4
$ws.Range("I3").Value = "x86-64"
$ws.Range("K3").Value = 32
$ws.Range("L3").Value = 2
$ws.Range("M3").Value = "DDR4"
$ws.Range("N3").Value = 3200
$ws.Range("O3").Value = 0.5
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.8
$ws.Range("R3").Value = 3.7

# Move the selection, matching the author's final cursor position
$ws.Range("S15").Select()
